# Update the Cd274-Cd80 LR-pairs sheet ("Natmi following Dr Hou advice"):
# recompute all existing data rows (sending x target cluster pairs for
# ECs/FAPs/M2) and add the previously-missing sCs sending-cluster rows,
# completing the full 4x4 cluster cross product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cd274"
$ws.Cells.Item(2,3).Value = "Cd80"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 12.230514
$ws.Cells.Item(2,8).Value = 36.691542
$ws.Cells.Item(2,9).Value = 0.4767524845277322
$ws.Cells.Item(2,10).Value = 0.4767524845277322
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.539733
$ws.Cells.Item(2,14).Value = 1.619199
$ws.Cells.Item(2,15).Value = 0.0338608224813917
$ws.Cells.Item(2,16).Value = 0.0338608224813917
$ws.Cells.Item(2,17).Value = 6.601212012762
$ws.Cells.Item(2,18).Value = 59.410908114858
$ws.Cells.Item(2,19).Value = 0.01614323124615598
$ws.Cells.Item(2,20).Value = 0.01614323124615598
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cd274"
$ws.Cells.Item(3,3).Value = "Cd80"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 12.230514
$ws.Cells.Item(3,8).Value = 36.691542
$ws.Cells.Item(3,9).Value = 0.4767524845277322
$ws.Cells.Item(3,10).Value = 0.4767524845277322
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.713472666666666
$ws.Cells.Item(3,14).Value = 11.140418
$ws.Cells.Item(3,15).Value = 0.2329693362375475
$ws.Cells.Item(3,16).Value = 0.2329693362375475
$ws.Cells.Item(3,17).Value = 45.41767943828399
$ws.Cells.Item(3,18).Value = 408.7591149445559
$ws.Cells.Item(3,19).Value = 0.1110687098700274
$ws.Cells.Item(3,20).Value = 0.1110687098700274
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Cd274"
$ws.Cells.Item(4,3).Value = "Cd80"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 12.230514
$ws.Cells.Item(4,8).Value = 36.691542
$ws.Cells.Item(4,9).Value = 0.4767524845277322
$ws.Cells.Item(4,10).Value = 0.4767524845277322
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 10.63672166666667
$ws.Cells.Item(4,14).Value = 31.910165
$ws.Cells.Item(4,15).Value = 0.6673079914308979
$ws.Cells.Item(4,16).Value = 0.6673079914308979
$ws.Cells.Item(4,17).Value = 130.09257325827
$ws.Cells.Item(4,18).Value = 1170.83315932443
$ws.Cells.Item(4,19).Value = 0.3181407428598912
$ws.Cells.Item(4,20).Value = 0.3181407428598912
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Cd274"
$ws.Cells.Item(5,3).Value = "Cd80"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 12.230514
$ws.Cells.Item(5,8).Value = 36.691542
$ws.Cells.Item(5,9).Value = 0.4767524845277322
$ws.Cells.Item(5,10).Value = 0.4767524845277322
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.049821333333333
$ws.Cells.Item(5,14).Value = 3.149464
$ws.Cells.Item(5,15).Value = 0.06586184985016284
$ws.Cells.Item(5,16).Value = 0.06586184985016284
$ws.Cells.Item(5,17).Value = 12.839854514832
$ws.Cells.Item(5,18).Value = 115.558690633488
$ws.Cells.Item(5,19).Value = 0.03139980055165758
$ws.Cells.Item(5,20).Value = 0.03139980055165759
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Cd274"
$ws.Cells.Item(6,3).Value = "Cd80"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.417539666666667
$ws.Cells.Item(6,8).Value = 7.252619
$ws.Cells.Item(6,9).Value = 0.09423708950643275
$ws.Cells.Item(6,10).Value = 0.09423708950643275
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.539733
$ws.Cells.Item(6,14).Value = 1.619199
$ws.Cells.Item(6,15).Value = 0.0338608224813917
$ws.Cells.Item(6,16).Value = 0.0338608224813917
$ws.Cells.Item(6,17).Value = 1.304825936909
$ws.Cells.Item(6,18).Value = 11.743433432181
$ws.Cells.Item(6,19).Value = 0.00319094535894034
$ws.Cells.Item(6,20).Value = 0.00319094535894034
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Cd274"
$ws.Cells.Item(7,3).Value = "Cd80"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.417539666666667
$ws.Cells.Item(7,8).Value = 7.252619
$ws.Cells.Item(7,9).Value = 0.09423708950643275
$ws.Cells.Item(7,10).Value = 0.09423708950643275
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.713472666666666
$ws.Cells.Item(7,14).Value = 11.140418
$ws.Cells.Item(7,15).Value = 0.2329693362375475
$ws.Cells.Item(7,16).Value = 0.2329693362375475
$ws.Cells.Item(7,17).Value = 8.97746747274911
$ws.Cells.Item(7,18).Value = 80.79720725474199
$ws.Cells.Item(7,19).Value = 0.02195435219127199
$ws.Cells.Item(7,20).Value = 0.02195435219127199
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Cd274"
$ws.Cells.Item(8,3).Value = "Cd80"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.417539666666667
$ws.Cells.Item(8,8).Value = 7.252619
$ws.Cells.Item(8,9).Value = 0.09423708950643275
$ws.Cells.Item(8,10).Value = 0.09423708950643275
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 10.63672166666667
$ws.Cells.Item(8,14).Value = 31.910165
$ws.Cells.Item(8,15).Value = 0.6673079914308979
$ws.Cells.Item(8,16).Value = 0.6673079914308979
$ws.Cells.Item(8,17).Value = 25.71469655245945
$ws.Cells.Item(8,18).Value = 231.432268972135
$ws.Cells.Item(8,19).Value = 0.06288516291683138
$ws.Cells.Item(8,20).Value = 0.06288516291683138
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Cd274"
$ws.Cells.Item(9,3).Value = "Cd80"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.417539666666667
$ws.Cells.Item(9,8).Value = 7.252619
$ws.Cells.Item(9,9).Value = 0.09423708950643275
$ws.Cells.Item(9,10).Value = 0.09423708950643275
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.049821333333333
$ws.Cells.Item(9,14).Value = 3.149464
$ws.Cells.Item(9,15).Value = 0.06586184985016284
$ws.Cells.Item(9,16).Value = 0.06586184985016284
$ws.Cells.Item(9,17).Value = 2.537984716246222
$ws.Cells.Item(9,18).Value = 22.841862446216
$ws.Cells.Item(9,19).Value = 0.00620662903938903
$ws.Cells.Item(9,20).Value = 0.00620662903938903
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Cd274"
$ws.Cells.Item(10,3).Value = "Cd80"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 10.06935333333333
$ws.Cells.Item(10,8).Value = 30.20806
$ws.Cells.Item(10,9).Value = 0.3925091962001163
$ws.Cells.Item(10,10).Value = 0.3925091962001163
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.539733
$ws.Cells.Item(10,14).Value = 1.619199
$ws.Cells.Item(10,15).Value = 0.0338608224813917
$ws.Cells.Item(10,16).Value = 0.0338608224813917
$ws.Cells.Item(10,17).Value = 5.43476228266
$ws.Cells.Item(10,18).Value = 48.91286054394
$ws.Cells.Item(10,19).Value = 0.01329068421484588
$ws.Cells.Item(10,20).Value = 0.01329068421484588
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Cd274"
$ws.Cells.Item(11,3).Value = "Cd80"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 10.06935333333333
$ws.Cells.Item(11,8).Value = 30.20806
$ws.Cells.Item(11,9).Value = 0.3925091962001163
$ws.Cells.Item(11,10).Value = 0.3925091962001163
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 3.713472666666666
$ws.Cells.Item(11,14).Value = 11.140418
$ws.Cells.Item(11,15).Value = 0.2329693362375475
$ws.Cells.Item(11,16).Value = 0.2329693362375475
$ws.Cells.Item(11,17).Value = 37.39226837434222
$ws.Cells.Item(11,18).Value = 336.5304153690799
$ws.Cells.Item(11,19).Value = 0.09144260690587439
$ws.Cells.Item(11,20).Value = 0.09144260690587439
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Cd274"
$ws.Cells.Item(12,3).Value = "Cd80"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 10.06935333333333
$ws.Cells.Item(12,8).Value = 30.20806
$ws.Cells.Item(12,9).Value = 0.3925091962001163
$ws.Cells.Item(12,10).Value = 0.3925091962001163
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 10.63672166666667
$ws.Cells.Item(12,14).Value = 31.910165
$ws.Cells.Item(12,15).Value = 0.6673079914308979
$ws.Cells.Item(12,16).Value = 0.6673079914308979
$ws.Cells.Item(12,17).Value = 107.1049087699889
$ws.Cells.Item(12,18).Value = 963.9441789299
$ws.Cells.Item(12,19).Value = 0.2619245233344558
$ws.Cells.Item(12,20).Value = 0.2619245233344558
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Cd274"
$ws.Cells.Item(13,3).Value = "Cd80"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 10.06935333333333
$ws.Cells.Item(13,8).Value = 30.20806
$ws.Cells.Item(13,9).Value = 0.3925091962001163
$ws.Cells.Item(13,10).Value = 0.3925091962001163
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.049821333333333
$ws.Cells.Item(13,14).Value = 3.149464
$ws.Cells.Item(13,15).Value = 0.06586184985016284
$ws.Cells.Item(13,16).Value = 0.06586184985016284
$ws.Cells.Item(13,17).Value = 10.57102194220444
$ws.Cells.Item(13,18).Value = 95.13919747984001
$ws.Cells.Item(13,19).Value = 0.02585138174494016
$ws.Cells.Item(13,20).Value = 0.02585138174494016
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Cd274"
$ws.Cells.Item(14,3).Value = "Cd80"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.9363953333333334
$ws.Cells.Item(14,8).Value = 2.809186
$ws.Cells.Item(14,9).Value = 0.03650122976571881
$ws.Cells.Item(14,10).Value = 0.03650122976571881
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.539733
$ws.Cells.Item(14,14).Value = 1.619199
$ws.Cells.Item(14,15).Value = 0.0338608224813917
$ws.Cells.Item(14,16).Value = 0.0338608224813917
$ws.Cells.Item(14,17).Value = 0.505403462446
$ws.Cells.Item(14,18).Value = 4.548631162014
$ws.Cells.Item(14,19).Value = 0.001235961661449495
$ws.Cells.Item(14,20).Value = 0.001235961661449495
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Cd274"
$ws.Cells.Item(15,3).Value = "Cd80"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.9363953333333334
$ws.Cells.Item(15,8).Value = 2.809186
$ws.Cells.Item(15,9).Value = 0.03650122976571881
$ws.Cells.Item(15,10).Value = 0.03650122976571881
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 3.713472666666666
$ws.Cells.Item(15,14).Value = 11.140418
$ws.Cells.Item(15,15).Value = 0.2329693362375475
$ws.Cells.Item(15,16).Value = 0.2329693362375475
$ws.Cells.Item(15,17).Value = 3.477278475527555
$ws.Cells.Item(15,18).Value = 31.295506279748
$ws.Cells.Item(15,19).Value = 0.008503667270373723
$ws.Cells.Item(15,20).Value = 0.008503667270373723
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Cd274"
$ws.Cells.Item(16,3).Value = "Cd80"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.9363953333333334
$ws.Cells.Item(16,8).Value = 2.809186
$ws.Cells.Item(16,9).Value = 0.03650122976571881
$ws.Cells.Item(16,10).Value = 0.03650122976571881
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 10.63672166666667
$ws.Cells.Item(16,14).Value = 31.910165
$ws.Cells.Item(16,15).Value = 0.6673079914308979
$ws.Cells.Item(16,16).Value = 0.6673079914308979
$ws.Cells.Item(16,17).Value = 9.960176530632221
$ws.Cells.Item(16,18).Value = 89.64158877569
$ws.Cells.Item(16,19).Value = 0.02435756231971952
$ws.Cells.Item(16,20).Value = 0.02435756231971952
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Cd274"
$ws.Cells.Item(17,3).Value = "Cd80"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 0.9363953333333334
$ws.Cells.Item(17,8).Value = 2.809186
$ws.Cells.Item(17,9).Value = 0.03650122976571881
$ws.Cells.Item(17,10).Value = 0.03650122976571881
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 1.049821333333333
$ws.Cells.Item(17,14).Value = 3.149464
$ws.Cells.Item(17,15).Value = 0.06586184985016284
$ws.Cells.Item(17,16).Value = 0.06586184985016284
$ws.Cells.Item(17,17).Value = 0.9830477973671111
$ws.Cells.Item(17,18).Value = 8.847430176304
$ws.Cells.Item(17,19).Value = 0.002404038514176067
$ws.Cells.Item(17,20).Value = 0.002404038514176067
